$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Cells.Item(9, 1).Value = 42726
$ws.Cells.Item(9, 2).Value = "6h"
$ws.Cells.Item(9, 3).Value = "FullCalendar Komponente eingebunden"

# Row 10
$ws.Cells.Item(10, 1).Value = 42727
$ws.Cells.Item(10, 2).Value = "10h"
$ws.Cells.Item(10, 3).Value = "Reservierungsvorgang erfolgreich implementiert"

# Copy date formatting from the previous date cell so the new cells reuse the
# same style (instead of creating a new number format / style entry).
$ws.Cells.Item(8, 1).Copy()
$ws.Range("A9:A10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("C11").Select()
